$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = "FA1_[FA-H]-"
$ws.Range("A3").Value = "FA2_[FA-H]-"
$ws.Range("A4").Value = "[LPL(FA1)-H]-"
$ws.Range("A5").Value = "[LPL(FA2)-H]-"
$ws.Range("A6").Value = "[LPL(FA1)-H2O-H]-"
$ws.Range("A7").Value = "[LPL(FA2)-H2O-H]-"

$ws.Range("A7").Select()
